# Logged Week 16 and performed season sim from Week 17
# Updates cumulative season stat sheets by appending the new week's
# per-drive/per-play logs (YDS, ST) and adding the week's totals to the
# running season totals (OFF, DEF, ST, TURNS, PEN).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet - append space-separated per-play yardage figures
# ---------------------------------------------------------------
$yds = $wb.Worksheets.Item("YDS")

$yds.Range("B2").Value = $yds.Range("B2").Value2 + " 3 3 -1 25 5 0 2 3 2 6 5 7 3 10 6 2 32 3 7 32 3 3 7 2"
$yds.Range("B3").Value = $yds.Range("B3").Value2 + " 7 41 0 9 8 9 13 13 26 3 3 24 9 -3 4 15"
$yds.Range("C2").Value = $yds.Range("C2").Value2 + " 4 11 2 2 -2 8 3 2 -2 0 0 1 -2 1 12 8 3 -4 4 4 -1 9 2 20 3 0 -1 1"
$yds.Range("C3").Value = $yds.Range("C3").Value2 + " 4 1 10 0 8 15 9 15 7 5 16 2 15 16 10 16 14 11 5 30 10 14 15"

# ---------------------------------------------------------------
# OFF sheet - add this week's offensive totals to season totals
# ---------------------------------------------------------------
$off = $wb.Worksheets.Item("OFF")

$off.Range("C2").Value = $off.Range("C2").Value2 + 15
$off.Range("D2").Value = $off.Range("D2").Value2 + 1
$off.Range("E2").Value = $off.Range("E2").Value2 + 1
$off.Range("F2").Value = $off.Range("F2").Value2 + 5
$off.Range("G2").Value = $off.Range("G2").Value2 + 2
$off.Range("L2").Value = $off.Range("L2").Value2 + 27
$off.Range("M2").Value = $off.Range("M2").Value2 + 16
$off.Range("O2").Value = $off.Range("O2").Value2 + 2
$off.Range("Q2").Value = $off.Range("Q2").Value2 + 59

$off.Range("C3").Value = $off.Range("C3").Value2 + 8
$off.Range("E3").Value = $off.Range("E3").Value2 + 2
$off.Range("F3").Value = $off.Range("F3").Value2 + 6
$off.Range("G3").Value = $off.Range("G3").Value2 + 2
$off.Range("H3").Value = $off.Range("H3").Value2 + 2
$off.Range("I3").Value = $off.Range("I3").Value2 + 3
$off.Range("J3").Value = $off.Range("J3").Value2 + 5
$off.Range("N3").Value = $off.Range("N3").Value2 + 2

# ---------------------------------------------------------------
# DEF sheet - add this week's defensive totals to season totals
# ---------------------------------------------------------------
$def = $wb.Worksheets.Item("DEF")

$def.Range("C2").Value = $def.Range("C2").Value2 + 12
$def.Range("D2").Value = $def.Range("D2").Value2 + 3
$def.Range("E2").Value = $def.Range("E2").Value2 + 1
$def.Range("F2").Value = $def.Range("F2").Value2 + 4
$def.Range("G2").Value = $def.Range("G2").Value2 + 5
$def.Range("J2").Value = $def.Range("J2").Value2 + 3
$def.Range("L2").Value = $def.Range("L2").Value2 + 35
$def.Range("M2").Value = $def.Range("M2").Value2 + 24
$def.Range("O2").Value = $def.Range("O2").Value2 + 3
$def.Range("P2").Value = $def.Range("P2").Value2 + 1
$def.Range("Q2").Value = $def.Range("Q2").Value2 + 73

$def.Range("B3").Value = $def.Range("B3").Value2 + 1
$def.Range("C3").Value = $def.Range("C3").Value2 + 14
$def.Range("E3").Value = $def.Range("E3").Value2 + 5
$def.Range("F3").Value = $def.Range("F3").Value2 + 5
$def.Range("G3").Value = $def.Range("G3").Value2 + 1
$def.Range("H3").Value = $def.Range("H3").Value2 + 5
$def.Range("I3").Value = $def.Range("I3").Value2 + 2
$def.Range("J3").Value = $def.Range("J3").Value2 + 4
$def.Range("N3").Value = $def.Range("N3").Value2 + 4

# ---------------------------------------------------------------
# ST sheet - special teams season totals + per-kick logs
# ---------------------------------------------------------------
$st = $wb.Worksheets.Item("ST")

$st.Range("B2").Value = $st.Range("B2").Value2 + 5
$st.Range("D2").Value = $st.Range("D2").Value2 + 4
$st.Range("F2").Value = $st.Range("F2").Value2 + 3
$st.Range("G2").Value = $st.Range("G2").Value2 + 3
$st.Range("J2").Value = $st.Range("J2").Value2 + 2
$st.Range("K2").Value = $st.Range("K2").Value2 + 1

$st.Range("B4").Value = $st.Range("B4").Value2 + " 60 62 60 53 54"
$st.Range("B5").Value = $st.Range("B5").Value2 + " 22 30 20 8 14"
$st.Range("B6").Value = $st.Range("B6").Value2 + " 23 40 13 16"
$st.Range("D3").Value = $st.Range("D3").Value2 + " 41 41 25 41"
$st.Range("D4").Value = $st.Range("D4").Value2 + " 4 28 0 0"
$st.Range("D5").Value = $st.Range("D5").Value2 + " 0 0 15"

# ---------------------------------------------------------------
# TURNS sheet - turnover season totals
# ---------------------------------------------------------------
$turns = $wb.Worksheets.Item("TURNS")

$turns.Range("E2").Value = $turns.Range("E2").Value2 + 1

# ---------------------------------------------------------------
# PEN sheet - penalty season totals
# ---------------------------------------------------------------
$pen = $wb.Worksheets.Item("PEN")

$pen.Range("B2").Value = $pen.Range("B2").Value2 + 2
$pen.Range("B3").Value = $pen.Range("B3").Value2 + 3
$pen.Range("D4").Value = $pen.Range("D4").Value2 + 2
